# ==========================================================================
# Edit script for ordenanzas/2028.docx
# Applies: keepNext + revised spacing across the document, bold VISTO/
# CONSIDERANDO captions (split off into their own paragraph), bold+indented
# "EL CONCEJO DELIBERANTE..." enactment clause (drop the leading "POR "),
# underlined "ARTICULO Nth" captions with the colon pulled into its own
# (still underlined) run, tightened "(" spacing in a few runs, and a
# pgNumType/start on the section.
# ==========================================================================

$d = $word.ActiveDocument

function Find-ParaByText($doc, $text) {
    $f = $doc.Content.Find
    $f.ClearFormatting()
    $ok = $f.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find-ParaByText: not found -> $text"
    }
    return $f.Parent.Paragraphs(1)
}

function Set-Spacing($para, $before, $after) {
    $para.Format.KeepWithNext = 1
    if ($null -ne $before) {
        $para.Format.SpaceBefore = $before
    }
    $para.Format.SpaceAfter = $after
}

# --------------------------------------------------------------------
# 1) Text-level fixes (do these while paragraph layout is still 1:1)
# --------------------------------------------------------------------

# Tighten the "(" runs that used to be padded with many spaces.
$f = $d.Content.Find
$f.ClearFormatting()
$f.Replacement.ClearFormatting()
$f.Execute("                    (", $true, $false, $false, $false, $false, $true, 1, $false, " (", 2) | Out-Null

# Drop the leading "POR " before the enactment clause.
$f = $d.Content.Find
$f.ClearFormatting()
$f.Replacement.ClearFormatting()
$f.Execute("POR EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA", $true, $false, $false, $false, $false, $true, 1, $false, "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA", 2) | Out-Null

Write-Output "text fixes done"

# --------------------------------------------------------------------
# 2) Split off "VISTO: " and "CONSIDERANDO: " into their own paragraphs,
#    each followed by a paragraph that starts with a plain space run.
# --------------------------------------------------------------------

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("VISTO: ", $true) | Out-Null
$visto = $f.Parent
$vistoIndex = $visto.Paragraphs(1).Index
$visto.InsertParagraphAfter()
$vistoEnd = $d.Paragraphs($vistoIndex).Range.End
$afterVisto = $d.Range($vistoEnd, $vistoEnd)
$afterVisto.InsertBefore(" ")

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("CONSIDERANDO: ", $true) | Out-Null
$considerando = $f.Parent
$considerandoIndex = $considerando.Paragraphs(1).Index
$considerando.InsertParagraphAfter()
$considerandoEnd = $d.Paragraphs($considerandoIndex).Range.End
$afterConsiderando = $d.Range($considerandoEnd, $considerandoEnd)
$afterConsiderando.InsertBefore(" ")

Write-Output "splits done"
Write-Output $d.Paragraphs.Count

# --------------------------------------------------------------------
# 3) Paragraph-level spacing / keepNext / alignment
# --------------------------------------------------------------------

# "Yerba Buena, 23 de Diciembre de 2015" — keepNext, after 200 -> 240
$p = Find-ParaByText $d "Yerba Buena"
Set-Spacing $p $null 12

# "ORDENANZA Nº 2028" — keepNext, before 240 / after 360, bold
$p = Find-ParaByText $d "ORDENANZA"
Set-Spacing $p 12 18
$p.Range.Bold = 1

# "VISTO: " caption — keepNext, before 240 / after 120, bold, drop justify
$p = Find-ParaByText $d "VISTO: "
Set-Spacing $p 12 6
$p.Format.Alignment = 0
$p.Range.Bold = 1

# "La naturaleza..." body — keepNext, after 120, drop justify
$p = Find-ParaByText $d "naturaleza Jurídica de los contratos"
Set-Spacing $p $null 6
$p.Format.Alignment = 0

# "CONSIDERANDO: " caption — keepNext, before 240 / after 120, bold, drop justify
$p = Find-ParaByText $d "CONSIDERANDO: "
Set-Spacing $p 12 6
$p.Format.Alignment = 0
$p.Range.Bold = 1

# "Que todo convenio..." body — keepNext, after 120, drop justify
$p = Find-ParaByText $d "Que todo convenio y/o contrato"
Set-Spacing $p $null 6
$p.Format.Alignment = 0

# "Que la relación contractual..." — keepNext, after 120, drop justify
$p = Find-ParaByText $d "Que la relación contractual mencionada"
Set-Spacing $p $null 6
$p.Format.Alignment = 0

# "Que es imprescindible..." — keepNext, after 120, drop justify
$p = Find-ParaByText $d "Que es imprescindible que dichos contratos"
Set-Spacing $p $null 6
$p.Format.Alignment = 0

# "EL CONCEJO DELIBERANTE SANCIONA..." — keepNext, before 360 / after 360,
# indent 1984/1984, bold (alignment already centered)
$p = Find-ParaByText $d "EL CONCEJO DELIBERANTE SANCIONA"
Set-Spacing $p 18 18
$p.Format.LeftIndent = 99.2
$p.Format.RightIndent = 99.2
$p.Range.Bold = 1

Write-Output "stage 3a done"

function Format-ArticuloHeading($doc, $heading) {
    # Underline the "ARTICULO Nth" caption word.
    $f = $doc.Content.Find
    $f.ClearFormatting()
    $f.Execute($heading, $true) | Out-Null
    $headRange = $f.Parent
    $headRange.Font.Underline = 1
    $para = $headRange.Paragraphs(1)

    # keepNext + after 120, drop justify
    $para.Format.KeepWithNext = 1
    $para.Format.SpaceAfter = 6
    $para.Format.Alignment = 0

    # Split the following ": " into an underlined ":" run and a plain " " run.
    $colon = $doc.Range($headRange.End, $headRange.End + 1)
    $colon.Font.Underline = 1
}

Format-ArticuloHeading $d "ARTICULO PRIMERO"
Format-ArticuloHeading $d "ARTICULO SEGUNDO"
Format-ArticuloHeading $d "ARTICULO TERCERO"
Format-ArticuloHeading $d "ARTICULO CUARTO"
Format-ArticuloHeading $d "ARTICULO QUINTO"
Format-ArticuloHeading $d "ARTICULO SEXTO"

Write-Output "articulo headings done"

# --------------------------------------------------------------------
# 4) List items under ARTICULO SEGUNDO — keepNext, after 200 -> 120,
#    drop justify (the numPr/ind stay untouched)
# --------------------------------------------------------------------

$listAnchors = @(
  "Fecha de realización.",
  "Dependencia/s municipal/es responsable.",
  "Funcionarios firmantes",
  "Norma legal o reglamentaria que",
  "Contraparte contractual",
  "Objeto del Convenio y/o Contrato.",
  "Detalle de la locación, prestación, servicio, obra o bien",
  "Erogación o acción que debe realizar la Municipalidad.",
  "Vigencia y monto total."
)
foreach ($anchor in $listAnchors) {
    $p = Find-ParaByText $d $anchor
    $p.Format.KeepWithNext = 1
    $p.Format.SpaceAfter = 6
    $p.Format.Alignment = 0
}

Write-Output "list items done"

# "El RECC contará además..." — keepNext, after 120, drop justify
$p = Find-ParaByText $d "El RECC contará además"
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# "ARTICULO SÉPTIMO..." — keepNext, after 120, drop justify (no run changes)
$p = Find-ParaByText $d "ARTICULO SÉPTIMO"
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# "ARTICULO OCTAVO..." — keepNext, after 120, drop justify (no run changes)
$p = Find-ParaByText $d "ARTICULO OCTAVO"
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

Write-Output "tail paragraphs done"

# --------------------------------------------------------------------
# 5) Section: restart page numbering at 3015
# --------------------------------------------------------------------

$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$ftr.PageNumbers.StartingNumber = 3015

Write-Output "page numbering done"
